$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalid-utilisation-report")

# Update the column header in H1 to reflect the new report wording
$ws.Range("H1").Value = "Fees paid to UKEF for the period"
